# Commit and Push as on 07-May-2019.
#
# - "Lipid Profile" was renamed to "Lipid profile" (lower-case "p") in both
#   the UploadLabReports and AddInvestigations sheets.
# - The active sheet moved from UploadLabReports (tab 1) to AddInvestigations
#   (tab 2), with the selection on that sheet moved to A4.

$wb = $excel.ActiveWorkbook

$wsUpload = $wb.Worksheets.Item("UploadLabReports")
$wsInvestigations = $wb.Worksheets.Item("AddInvestigations")

# Fix the investigation-type label's capitalisation on both sheets that
# reference it.
$wsUpload.Range("A4").Value = "Lipid profile"
$wsInvestigations.Range("A4").Value = "Lipid profile"

# Make "AddInvestigations" the active sheet/tab, with A4 selected.
$wsInvestigations.Activate()
$wsInvestigations.Range("A4").Select()
